$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Zi Hao Tan (row 8) - update back end function lists role description
$ws.Range("C8").Value = "Created a template for back end function lists and edit. Created the userlist page (P9) `nand edit an order profile page (P12). Created the meat aisle (P2) and meat product `ndescription page (P3). Created the CSS related to those pages."

# James Partsafas (row 2) - update Homepage role description
$ws.Range("C2").Value = "Homepage (P1) (index.html) and everything therein, including the banners, food displays, buttons, and basic styles that are applied elsewhere. Also made (P2) grain aisle page and all (P3) pages relating to grain foods. I made the user edit page (P10). I also handled hosting of the website and I was the Github administrator, responsible for handling merge conflicts and other issues."

# Update selection to C2
$ws.Range("C2").Select()
